# Remove the "Fuel Technology Uranium ELC" (FT-ELCURN) entry.
#
# This entry occupies row 17 (the ~FI_T / Fuel Technology table row for
# FT-ELCURN) and row 33 (the \I: Process Set Membership row that links
# back to it via =B17 / =C17 formulas). Deleting both rows shifts every
# row below each deletion point up by one, which matches the target
# layout exactly (comments move from row 21->20, 29->28, 30->29, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the higher-numbered row first so the row 17 index used below is
# still valid when we get to it.
$ws.Rows.Item(33).Delete()
$ws.Rows.Item(17).Delete()

$ws.Range("N21").Select()
